$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2063.3572
$ws.Range("I19").Value = 2933.375
$ws.Range("K19").Value = 2933.375
$ws.Range("M19").Value = -2758.375
$ws.Range("H20").Value = 10521
$ws.Range("I20").Value = 10521
$ws.Range("K20").Value = 10521
$ws.Range("M20").Value = -10291
$ws.Range("H35").Value = 10521
$ws.Range("I35").Value = 10521
$ws.Range("K35").Value = 10521
$ws.Range("M35").Value = -10142
$ws.Range("H43").Value = 1194.0714
$ws.Range("I43").Value = 987.75
$ws.Range("J43").Value = 1276.6
$ws.Range("K43").Value = 987.75
$ws.Range("L43").Value = 1276.6
$ws.Range("M43").Value = -918.75
$ws.Range("N43").Value = -1414.6
$ws.Range("H74").Value = 15629563
$ws.Range("I74").Value = 4167.6665
$ws.Range("J74").Value = 25004800
$ws.Range("K74").Value = 4167.6665
$ws.Range("L74").Value = 25004800
$ws.Range("M74").Value = -3231.6665
$ws.Range("N74").Value = -25006672
$ws.Range("H77").Value = 15629563
$ws.Range("I77").Value = 4167.6665
$ws.Range("J77").Value = 25004800
$ws.Range("K77").Value = 20838.3325
$ws.Range("L77").Value = 125024000
$ws.Range("M77").Value = -16158.3325
$ws.Range("N77").Value = -125033360
$ws.Range("H112").Value = 2924948
$ws.Range("J112").Value = 3087424.2
$ws.Range("L112").Value = 9262272.600000001
$ws.Range("N112").Value = -9264488.600000001
$ws.Range("H129").Value = 846.9194
$ws.Range("J129").Value = 849.47455
$ws.Range("L129").Value = 2548.42365
$ws.Range("N129").Value = -12548.42365
$ws.Range("H132").Value = 4067.6316
$ws.Range("I132").Value = 4486.5625
$ws.Range("K132").Value = 13459.6875
$ws.Range("M132").Value = -10929.6875
$ws.Range("H141").Value = 2342
$ws.Range("I141").Value = 1690
$ws.Range("J141").Value = 4950
$ws.Range("K141").Value = 5070
$ws.Range("L141").Value = 14850
$ws.Range("M141").Value = 110
$ws.Range("N141").Value = -25210

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3204.1555
$ws.Range("I45").Value = 2937.389
$ws.Range("K45").Value = 2937.389
$ws.Range("M45").Value = -2560.389
$ws.Range("H61").Value = 4566.7334
$ws.Range("I61").Value = 3789
$ws.Range("K61").Value = 3789
$ws.Range("M61").Value = -3577
$ws.Range("H74").Value = 1937.25
$ws.Range("I74").Value = 2251.7646
$ws.Range("J74").Value = 1173.4286
$ws.Range("K74").Value = 2251.7646
$ws.Range("L74").Value = 1173.4286
$ws.Range("M74").Value = -1377.7646
$ws.Range("N74").Value = -2921.4286
$ws.Range("H77").Value = 1937.25
$ws.Range("I77").Value = 2251.7646
$ws.Range("J77").Value = 1173.4286
$ws.Range("K77").Value = 11258.823
$ws.Range("L77").Value = 5867.143
$ws.Range("M77").Value = -6890.823
$ws.Range("N77").Value = -14603.143
$ws.Range("H97").Value = 1602.5883
$ws.Range("I97").Value = 1450.9333
$ws.Range("J97").Value = 2740
$ws.Range("K97").Value = 1450.9333
$ws.Range("L97").Value = 2740
$ws.Range("M97").Value = -954.9332999999999
$ws.Range("N97").Value = -3732
$ws.Range("H132").Value = 19854.379
$ws.Range("I132").Value = 2322.9285
$ws.Range("K132").Value = 6968.7855
$ws.Range("M132").Value = -4438.7855
$ws.Range("H136").Value = 4566.7334
$ws.Range("I136").Value = 3789
$ws.Range("K136").Value = 11367
$ws.Range("M136").Value = -8817
$ws.Range("H138").Value = 33821.75
$ws.Range("J138").Value = 33821.75
$ws.Range("L138").Value = 33821.75
$ws.Range("N138").Value = -44101.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1200
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 7000
$ws.Range("J4").Value = 7000
$ws.Range("L4").Value = 7000
$ws.Range("N4").Value = -7224
$ws.Range("H58").Value = 19143.607
$ws.Range("I58").Value = 1516.7333
$ws.Range("J58").Value = 39482.31
$ws.Range("K58").Value = 1516.7333
$ws.Range("L58").Value = 39482.31
$ws.Range("M58").Value = -1313.7333
$ws.Range("N58").Value = -39888.31
$ws.Range("H107").Value = 2064.125
$ws.Range("J107").Value = 3431.5
$ws.Range("L107").Value = 3431.5
$ws.Range("N107").Value = -7271.5
$ws.Range("H114").Value = 35842.5
$ws.Range("J114").Value = 35842.5
$ws.Range("L114").Value = 35842.5
$ws.Range("N114").Value = -44520.5
$ws.Range("H132").Value = 14727.875
$ws.Range("I132").Value = 17742.258
$ws.Range("J132").Value = 4345
$ws.Range("K132").Value = 53226.774
$ws.Range("L132").Value = 13035
$ws.Range("M132").Value = -50696.774
$ws.Range("N132").Value = -18095
$ws.Range("H136").Value = 19143.607
$ws.Range("I136").Value = 1516.7333
$ws.Range("J136").Value = 39482.31
$ws.Range("K136").Value = 4550.199900000001
$ws.Range("L136").Value = 118446.93
$ws.Range("M136").Value = -2000.199900000001
$ws.Range("N136").Value = -123546.93

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("H5").Value = 1846.2727
$ws.Range("J5").Value = 3002.5
$ws.Range("L5").Value = 9007.5
$ws.Range("N5").Value = -9231.5
$ws.Range("H69").Value = 1350
$ws.Range("I69").Value = 200
$ws.Range("K69").Value = 600
$ws.Range("M69").Value = 211
$ws.Range("H72").Value = 1350
$ws.Range("I72").Value = 200
$ws.Range("K72").Value = 1800
$ws.Range("M72").Value = 2256
$ws.Range("H107").Value = 20314.8
$ws.Range("J107").Value = 393.5
$ws.Range("L107").Value = 1180.5
$ws.Range("N107").Value = -5020.5
$ws.Range("H131").Value = 749.28
$ws.Range("J131").Value = 765.1875
$ws.Range("L131").Value = 2295.5625
$ws.Range("N131").Value = -12375.5625
$ws.Range("H135").Value = 1846.2727
$ws.Range("J135").Value = 3002.5
$ws.Range("L135").Value = 27022.5
$ws.Range("N135").Value = -32092.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H53").Value = 3815.6
$ws.Range("I53").Value = 3519.5
$ws.Range("J53").Value = 5000
$ws.Range("K53").Value = 3519.5
$ws.Range("L53").Value = 5000
$ws.Range("M53").Value = -2888.5
$ws.Range("N53").Value = -6262
$ws.Range("H97").Value = 1725.3793
$ws.Range("I97").Value = 1170.5238
$ws.Range("J97").Value = 3181.875
$ws.Range("K97").Value = 1170.5238
$ws.Range("L97").Value = 3181.875
$ws.Range("M97").Value = -674.5237999999999
$ws.Range("N97").Value = -4173.875
$ws.Range("H107").Value = 4150.75
$ws.Range("H108").Value = 45331.668
$ws.Range("J108").Value = 45331.668
$ws.Range("L108").Value = 45331.668
$ws.Range("N108").Value = -53011.668
$ws.Range("H117").Value = 20000
$ws.Range("J117").Value = 20000
$ws.Range("L117").Value = 20000
$ws.Range("N117").Value = -26884
$ws.Range("H132").Value = 78145.10000000001
$ws.Range("I132").Value = 74314.57000000001
$ws.Range("K132").Value = 222943.71
$ws.Range("M132").Value = -220413.71
$ws.Range("H138").Value = 55000
$ws.Range("J138").Value = 55000
$ws.Range("L138").Value = 55000
$ws.Range("N138").Value = -65280

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1685000
$ws.Range("I2").Value = 1857142.9
$ws.Range("J2").Value = 480000
$ws.Range("K2").Value = 1857142.9
$ws.Range("L2").Value = 480000
$ws.Range("M2").Value = -1857030.9
$ws.Range("N2").Value = -480224
$ws.Range("H16").Value = 612.5
$ws.Range("I16").Value = 675
$ws.Range("J16").Value = 550
$ws.Range("K16").Value = 675
$ws.Range("L16").Value = 550
$ws.Range("M16").Value = -505
$ws.Range("N16").Value = -890
$ws.Range("H61").Value = 4915.75
$ws.Range("I61").Value = 2043.6364
$ws.Range("K61").Value = 2043.6364
$ws.Range("M61").Value = -1841.6364
$ws.Range("H113").Value = 4915.75
$ws.Range("I113").Value = 2043.6364
$ws.Range("K113").Value = 2043.6364
$ws.Range("M113").Value = 126.3635999999999
$ws.Range("H132").Value = 1990
$ws.Range("I132").Value = 1354.75
$ws.Range("K132").Value = 4064.25
$ws.Range("M132").Value = -1534.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 40001.5
$ws.Range("J2").Value = 70003
$ws.Range("L2").Value = 70003
$ws.Range("N2").Value = -70227
$ws.Range("H107").Value = 4547088.5
$ws.Range("I107").Value = 584.75
$ws.Range("J107").Value = 7578090.5
$ws.Range("K107").Value = 1754.25
$ws.Range("L107").Value = 22734271.5
$ws.Range("M107").Value = 165.75
$ws.Range("N107").Value = -22738111.5
$ws.Range("H132").Value = 3141.84
$ws.Range("I132").Value = 2987.2222
$ws.Range("J132").Value = 3539.4285
$ws.Range("K132").Value = 8961.6666
$ws.Range("L132").Value = 10618.2855
$ws.Range("M132").Value = -6431.6666
$ws.Range("N132").Value = -15678.2855
$ws.Range("H141").Value = 56333
$ws.Range("J141").Value = 56333
$ws.Range("L141").Value = 56333
$ws.Range("N141").Value = -66693
